$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4903.377669225691
$ws.Range("C2").Value = 20525.24663867073
$ws.Range("D2").Value = 19482.449132881
